$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 251.77777
$ws.Range("I5").Value = 64.14286
$ws.Range("J5").Value = 371.18182
$ws.Range("K5").Value = 64.14286
$ws.Range("L5").Value = 371.18182
$ws.Range("M5").Value = 50.85714
$ws.Range("N5").Value = -601.18182
$ws.Range("H15").Value = 264.91
$ws.Range("I15").Value = 264.91
$ws.Range("K15").Value = 794.73
$ws.Range("M15").Value = -625.73
$ws.Range("H18").Value = 714.25
$ws.Range("I18").Value = 597.36365
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 597.36365
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -313.36365
$ws.Range("N18").Value = -2568
$ws.Range("H33").Value = 4330240.5
$ws.Range("I33").Value = 1547.6666
$ws.Range("J33").Value = 15151972
$ws.Range("K33").Value = 1547.6666
$ws.Range("L33").Value = 15151972
$ws.Range("M33").Value = -1318.6666
$ws.Range("N33").Value = -15152430
$ws.Range("H113").Value = 8335796
$ws.Range("I113").Value = 11113077
$ws.Range("J113").Value = 3953
$ws.Range("K113").Value = 11113077
$ws.Range("L113").Value = 3953
$ws.Range("M113").Value = -11109823
$ws.Range("N113").Value = -10461
$ws.Range("H132").Value = 1785.6222
$ws.Range("I132").Value = 1487.5581
$ws.Range("K132").Value = 4462.6743
$ws.Range("M132").Value = -1932.6743
$ws.Range("H138").Value = 2073.5898
$ws.Range("I138").Value = 1332.037
$ws.Range("J138").Value = 2466.1765
$ws.Range("K138").Value = 3996.111
$ws.Range("L138").Value = 7398.529500000001
$ws.Range("M138").Value = 1143.889
$ws.Range("N138").Value = -17678.5295

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3615.97
$ws.Range("I32").Value = 3112.5852
$ws.Range("J32").Value = 11502.333
$ws.Range("K32").Value = 3112.5852
$ws.Range("L32").Value = 11502.333
$ws.Range("M32").Value = -2825.5852
$ws.Range("N32").Value = -12076.333
$ws.Range("H61").Value = 248636.73
$ws.Range("I61").Value = 5681.3335
$ws.Range("J61").Value = 591632.5600000001
$ws.Range("K61").Value = 5681.3335
$ws.Range("L61").Value = 591632.5600000001
$ws.Range("M61").Value = -5469.3335
$ws.Range("N61").Value = -592056.5600000001
$ws.Range("H132").Value = 1565120.2
$ws.Range("I132").Value = 1684.5306
$ws.Range("J132").Value = 6672343.5
$ws.Range("K132").Value = 5053.5918
$ws.Range("L132").Value = 20017030.5
$ws.Range("M132").Value = -2523.5918
$ws.Range("N132").Value = -20022090.5
$ws.Range("H136").Value = 248636.73
$ws.Range("I136").Value = 5681.3335
$ws.Range("J136").Value = 591632.5600000001
$ws.Range("K136").Value = 17044.0005
$ws.Range("L136").Value = 1774897.68
$ws.Range("M136").Value = -14494.0005
$ws.Range("N136").Value = -1779997.68

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 24861.459
$ws.Range("I134").Value = 4737
$ws.Range("J134").Value = 69135.266
$ws.Range("K134").Value = 14211
$ws.Range("L134").Value = 207405.798
$ws.Range("M134").Value = -11676
$ws.Range("N134").Value = -212475.798

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2051.7368
$ws.Range("I16").Value = 1892.6666
$ws.Range("J16").Value = 2648.25
$ws.Range("K16").Value = 1892.6666
$ws.Range("L16").Value = 2648.25
$ws.Range("M16").Value = -1605.6666
$ws.Range("N16").Value = -3222.25
$ws.Range("H31").Value = 4200.278
$ws.Range("I31").Value = 1927.7906
$ws.Range("J31").Value = 7569.8276
$ws.Range("K31").Value = 1927.7906
$ws.Range("L31").Value = 7569.8276
$ws.Range("M31").Value = -1632.7906
$ws.Range("N31").Value = -8159.8276
$ws.Range("H34").Value = 4200.278
$ws.Range("I34").Value = 1927.7906
$ws.Range("J34").Value = 7569.8276
$ws.Range("K34").Value = 1927.7906
$ws.Range("L34").Value = 7569.8276
$ws.Range("M34").Value = -1725.7906
$ws.Range("N34").Value = -7973.8276
$ws.Range("H107").Value = 637
$ws.Range("I107").Value = 501.6842
$ws.Range("J107").Value = 820.6429000000001
$ws.Range("K107").Value = 501.6842
$ws.Range("L107").Value = 820.6429000000001
$ws.Range("M107").Value = 1418.3158
$ws.Range("N107").Value = -4660.6429
$ws.Range("H113").Value = 2051.7368
$ws.Range("I113").Value = 1892.6666
$ws.Range("J113").Value = 2648.25
$ws.Range("K113").Value = 1892.6666
$ws.Range("L113").Value = 2648.25
$ws.Range("M113").Value = 277.3334
$ws.Range("N113").Value = -6988.25
$ws.Range("H122").Value = 1390967.2
$ws.Range("I122").Value = 1853783
$ws.Range("J122").Value = 2520
$ws.Range("K122").Value = 5561349
$ws.Range("L122").Value = 7560
$ws.Range("M122").Value = -5558899
$ws.Range("N122").Value = -12460
$ws.Range("H134").Value = 198485.05
$ws.Range("I134").Value = 2330.6304
$ws.Range("K134").Value = 6991.8912
$ws.Range("M134").Value = -4456.8912

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1409.1714
$ws.Range("I5").Value = 456.1
$ws.Range("J5").Value = 1790.4
$ws.Range("K5").Value = 1368.3
$ws.Range("L5").Value = 5371.200000000001
$ws.Range("M5").Value = -1256.3
$ws.Range("N5").Value = -5595.200000000001
$ws.Range("H114").Value = 4089.4285
$ws.Range("I114").Value = 250.41667
$ws.Range("J114").Value = 6092.391
$ws.Range("K114").Value = 751.25001
$ws.Range("L114").Value = 18277.173
$ws.Range("M114").Value = 2502.74999
$ws.Range("N114").Value = -24785.173
$ws.Range("H135").Value = 1409.1714
$ws.Range("I135").Value = 456.1
$ws.Range("J135").Value = 1790.4
$ws.Range("K135").Value = 4104.900000000001
$ws.Range("L135").Value = 16113.6
$ws.Range("M135").Value = -1569.900000000001
$ws.Range("N135").Value = -21183.6

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 43829080
$ws.Range("I122").Value = 46298010
$ws.Range("K122").Value = 138894030
$ws.Range("M122").Value = -138891580
$ws.Range("H132").Value = 3586.6667
$ws.Range("I132").Value = 4068.75
$ws.Range("K132").Value = 12206.25
$ws.Range("M132").Value = -9676.25

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1461.6428
$ws.Range("I61").Value = 1460.3636
$ws.Range("J61").Value = 1466.3334
$ws.Range("K61").Value = 1460.3636
$ws.Range("L61").Value = 1466.3334
$ws.Range("M61").Value = -1258.3636
$ws.Range("N61").Value = -1870.3334
$ws.Range("H64").Value = 24320
$ws.Range("J64").Value = 24320
$ws.Range("L64").Value = 24320
$ws.Range("N64").Value = -24770
$ws.Range("H67").Value = 24320
$ws.Range("J67").Value = 24320
$ws.Range("L67").Value = 24320
$ws.Range("N67").Value = -25880
$ws.Range("H113").Value = 1461.6428
$ws.Range("I113").Value = 1460.3636
$ws.Range("J113").Value = 1466.3334
$ws.Range("K113").Value = 1460.3636
$ws.Range("L113").Value = 1466.3334
$ws.Range("M113").Value = 709.6364000000001
$ws.Range("N113").Value = -5806.3334

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 41562.25
$ws.Range("J63").Value = 41562.25
$ws.Range("L63").Value = 41562.25
$ws.Range("N63").Value = -42810.25
$ws.Range("H66").Value = 41562.25
$ws.Range("J66").Value = 41562.25
$ws.Range("L66").Value = 124686.75
$ws.Range("N66").Value = -130926.75
$ws.Range("H100").Value = 77285.84
$ws.Range("I100").Value = 91156.91
$ws.Range("K100").Value = 182313.82
$ws.Range("M100").Value = -181772.82
$ws.Range("H126").Value = 661.4074000000001
$ws.Range("I126").Value = 561.381
$ws.Range("J126").Value = 1011.5
$ws.Range("K126").Value = 1684.143
$ws.Range("L126").Value = 3034.5
$ws.Range("M126").Value = 785.857
$ws.Range("N126").Value = -7974.5
$ws.Range("H136").Value = 2003.1587
$ws.Range("I136").Value = 1889.1428
$ws.Range("J136").Value = 2231.1904
$ws.Range("K136").Value = 5667.428400000001
$ws.Range("L136").Value = 6693.5712
$ws.Range("M136").Value = -3117.428400000001
$ws.Range("N136").Value = -11793.5712
